$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaIns = $metaPara.Range
$metaIns.Collapse(1)
$metaIns.InsertAfter("Meta description: Read our review of Big Thunder King Strike, an online slot game from Ainsworth with free spins, jackpots, and jungle-themed symbols. Play for free now.")

# Make just the "Meta description" lead-in bold (leave the rest plain).
$metaStart = $d.Paragraphs(2).Range.Start
$metaBoldEnd = $metaStart + "Meta description".Length
$metaBoldRange = $d.Range($metaStart, $metaBoldEnd)
$metaBoldRange.Bold = 1

# ------------------------------------------------------------------
# 2. Remove the duplicated bold title paragraph
#    ("Play Big Thunder King Strike Free Slot | Ainsworth Gaming")
#    that used to sit right before the closing italic paragraph.
# ------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
for ($i = $paraCount; $i -ge 1; $i--) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -eq "Play Big Thunder King Strike Free Slot | Ainsworth Gaming" `
        -and $candidate.Range.Font.Bold) {
        $candidate.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 3. Replace the text of the final (italic) paragraph with the new
#    image-prompt copy, keeping its italic formatting intact.
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastStart = $d.Paragraphs($lastIndex).Range.Start
$lastEnd = $d.Paragraphs($lastIndex).Range.End
$lastRange = $d.Range($lastStart, $lastEnd)
$lastRange.Text = "Please create a feature image for ""Big Thunder King Strike"" that fits the theme of the game and features a happy Maya warrior wearing glasses. The image should be in a cartoon style and can include elements such as jungle foliage, animals, or tribal weapons. Be creative and use bold colors to capture the game's eccentric and adventurous vibe."
